# Weekly price-sheet update: a new week's price record (2023-09-06) is
# inserted as a new data row right after the header/first data block at
# row 50, pushing all the existing historical rows (old 50..118) down by
# one (new 51..119). The sheet's used range grows from A1:R118 to
# A1:R119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 50; rows 50-118 shift down to 51-119.
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the new week's data.
$ws.Cells.Item(50, 1).Value  = 9
$ws.Cells.Item(50, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(50, 3).Value  = "Metropolitana"
$ws.Cells.Item(50, 4).Value  = 45175
$ws.Cells.Item(50, 5).Value  = 13
$ws.Cells.Item(50, 6).Value  = 100112035
$ws.Cells.Item(50, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(50, 8).Value  = "Sin especificar"
$ws.Cells.Item(50, 9).Value  = "Primera"
$ws.Cells.Item(50, 10).Value = 52
$ws.Cells.Item(50, 11).Value = 19000
$ws.Cells.Item(50, 12).Value = 21000
$ws.Cells.Item(50, 13).Value = 20000
$ws.Cells.Item(50, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(50, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(50, 16).Value = 1333
$ws.Cells.Item(50, 17).Value = 15
$ws.Cells.Item(50, 18).Value = "Hortaliza"
